# Fruta / hortaliza, semanal
# Update the data rows so that the weekly records for each market visit are
# re-ordered / updated to reflect the latest source data (columns D, J, K,
# L, M, P on rows 2 and 4-10; row 3 is unchanged).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = @{ D = 45091; J = 40;  K = 20000; L = 22000; M = 21000; P = 1400 }
    4  = @{ D = 45133; J = 50;  K = 22000; L = 22000; M = 22000; P = 1467 }
    5  = @{ D = 44750; J = 140; K = 19000; L = 20000; M = 19571; P = 1305 }
    6  = @{ D = 45084; J = 90;  K = 22000; L = 23000; M = 22556; P = 1504 }
    7  = @{ D = 45119; J = 50;  K = 20000; L = 20000; M = 20000; P = 1333 }
    8  = @{ D = 45141; J = 50;  K = 8500;  L = 9000;  M = 8800;  P = 587 }
    9  = @{ D = 44749; J = 90;  K = 17000; L = 18000; M = 17556; P = 1170 }
    10 = @{ D = 45063; J = 40;  K = 21000; L = 22000; M = 21500; P = 1433 }
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("J$row").Value = $vals.J
    $ws.Range("K$row").Value = $vals.K
    $ws.Range("L$row").Value = $vals.L
    $ws.Range("M$row").Value = $vals.M
    $ws.Range("P$row").Value = $vals.P
}
